# Rename the original (only) sheet from "Sheet1" to "Before"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Before"

# Add a new worksheet right after "Before" and name it "New"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "New"

# Header row
$ws2.Cells.Item(1, 1).Value = "idcity"
$ws2.Cells.Item(1, 2).Value = "ncity"

# Data: idcity values (column A) and their matching ncity group (column B),
# in the same order as the "New" sheet in the target workbook.
$ids = @(1101,1401,1501,1601,1701,1801,2101,2201,2301,2401,2501,2801,2901,3201,3801,5001,301,302,601,602,1901,1902,2601,2602,3701,3702,3901,3902,4701,4702,4801,4802,101,102,103,3301,3302,3303,3401,3402,3403,4301,4302,4303,1001,1002,1003,1004,3601,3602,3603,3604,3605,3606,501,502,503,504,505,506,507,4401,4402,4403,4404,4405,4406,4407)
$ns  = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,3,3,3,4,4,4,4,6,6,6,6,6,6,7,7,7,7,7,7,7,7,7,7,7,7,7,7)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws2.Cells.Item($i + 2, 1).Value = $ids[$i]
    $ws2.Cells.Item($i + 2, 2).Value = $ns[$i]
}

# Match the author's final column width / selection state on the new sheet
$ws2.Columns.Item(1).ColumnWidth = 9.1666666666667
$selResult = $ws2.Range("B2").Select()
